# Applies the OOXML changes described by the diff:
#   - adds <w:pStyle w:val="Normal"/> as the first child of <w:pPr> in the
#     eight affected paragraphs ("Au fond de mon coeur..." through
#     "Le chat est orange.")
#   - splits the single <w:r> that carries the paragraph's text into two
#     runs: an empty leading run (<w:r><w:rPr/></w:r>) followed by a run
#     that holds the text, both with an empty <w:rPr/> (so the <w:lang>
#     run formatting, any w:rsidRPr attribute, and any
#     <w:lastRenderedPageBreak/> marker are dropped).
#
# The transform is driven off each paragraph's own current OOXML (fetched
# via Range.WordOpenXML) so paragraph-level identity attributes
# (w14:paraId, w:rsidR, ...) and any other <w:pPr> children (w:ind, w:jc,
# the paragraph mark's w:rPr, ...) are preserved byte-for-byte; only the
# <w:pStyle> insertion and the run split are performed.

$d = $word.ActiveDocument

function Get-FixedParagraphXml($para) {
    $full = $para.Range.WordOpenXML()

    $partRe = [regex]'(?s)<pkg:part pkg:name="/word/document\.xml"[^>]*><pkg:xmlData>(.*?)</pkg:xmlData></pkg:part>'
    $partMatch = $partRe.Match($full)
    $docXml = $partMatch.Groups[1].Value

    $bodyRe = [regex]'(?s)<w:body>(.*?)</w:body>'
    $bodyMatch = $bodyRe.Match($docXml)
    $bodyInner = $bodyMatch.Groups[1].Value

    $paraRe = [regex]'(?s)^<w:p\b.*?</w:p>'
    $paraMatch = $paraRe.Match($bodyInner)
    $paraXml = $paraMatch.Value

    # 1) Add <w:pStyle w:val="Normal"/> as the first child of <w:pPr>.
    $paraXml = [regex]::Replace($paraXml, '<w:pPr>', '<w:pPr><w:pStyle w:val="Normal"/>', 1)

    # 2) Gather the run text (concatenation of every <w:t> in the paragraph).
    $tRe = [regex]'(?s)<w:t[^>]*>(.*?)</w:t>'
    $sb = New-Object System.Text.StringBuilder
    foreach ($tm in $tRe.Matches($paraXml)) {
        [void]$sb.Append($tm.Groups[1].Value)
    }
    $text = $sb.ToString()

    # 3) Drop every existing run (this also removes stray run-level
    #    markers such as <w:lastRenderedPageBreak/>).
    $paraXml = [regex]::Replace($paraXml, '(?s)<w:r\b.*?</w:r>', '')
    $paraXml = [regex]::Replace($paraXml, '(?s)<w:r\b[^>]*/>', '')

    # 4) Re-insert the runs as an empty run followed by a text run, both
    #    with an empty <w:rPr/>.
    $newRuns = '<w:r><w:rPr/></w:r><w:r><w:rPr/><w:t>' + $text + '</w:t></w:r>'
    $paraXml = $paraXml -replace '</w:p>$', ($newRuns + '</w:p>')

    $pkg = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
           '<w:body>' + $paraXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    return $pkg
}

# The eight paragraphs targeted by the diff, identified by a distinctive
# text fragment so the script is resilient to paragraph-index drift.
$targets = @(
    "Au fond de mon c",
    "Compar",
    "J'aurais aim",
    "Sauver une fille",
    "Attendez une minute.",
    "Un jour, un myst",
    "Ou que dites-vous de",
    "Le chat est orange."
)

foreach ($needle in $targets) {
    $found = $null
    foreach ($para in $d.Paragraphs) {
        if ($para.Range.Text.StartsWith($needle)) {
            $found = $para
            break
        }
    }
    if ($found -eq $null) {
        Write-Output "NOT FOUND: $needle"
        continue
    }
    $xml = Get-FixedParagraphXml $found
    [void]$found.Range.InsertXML($xml)
    Write-Output "fixed: $needle"
}
